$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '87.446.15'
$ws.Range("E2").Value = '  -0.76%  '
$ws.Range("D3").Value = '3.217.93'
$ws.Range("E3").Value = '  -3.59%  '
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '204.15'
$ws.Range("E5").Value = '  -7.35%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '609.66'
$ws.Range("E6").Value = '  -6.59%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.377'
$ws.Range("E7").Value = '  +3.45%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.674'
$ws.Range("E8").Value = '  +10.96%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.998'
$ws.Range("E9").Value = '  -0.09%  '
$ws.Range("D10").Value = '3.213.37'
$ws.Range("E10").Value = '  -3.62%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.540'
$ws.Range("E11").Value = '  -8.45%  '
$ws.Range("E12").Value = '  +6.42%  '
$ws.Range("E13").Value = '  -9.47%  '
$ws.Range("D14").Value = '3.807.80'
$ws.Range("E14").Value = '  -3.59%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '5.29'
$ws.Range("E15").Value = '  -4.27%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '32.57'
$ws.Range("E16").Value = '  -9.26%  '
$ws.Range("D17").Value = '87.205.98'
$ws.Range("E17").Value = '  -0.77%  '
$ws.Range("D18").Value = '3.221.60'
$ws.Range("E18").Value = '  -3.18%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '13.50'
$ws.Range("E19").Value = '  -8.59%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '2.95'
$ws.Range("E20").Value = '  -6.49%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '421.05'
$ws.Range("E21").Value = '  -8.31%  '
$ws.Range("E22").Value = '  -12.92%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.12'
$ws.Range("E23").Value = '  -7.82%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '5.20'
$ws.Range("E24").Value = '  -8.69%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '11.61'
$ws.Range("E25").Value = '  -9.93%  '
$ws.Range("D26").Value = '3.386.90'
$ws.Range("E26").Value = '  -3.37%  '
$ws.Range("B27").Value = 'Litecoin'
$ws.Range("C27").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '74.19'
$ws.Range("E27").Value = '  -6.20%  '
$ws.Range("B28").Value = 'PEPE'
$ws.Range("C28").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.0000132'
$ws.Range("E28").Value = '  +3.70%  '
$ws.Range("E29").Value = '  -0.02%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.171'
$ws.Range("E30").Value = '  -15.30%  '
$ws.Range("E31").Value = '  +0.04%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '545.95'
$ws.Range("E32").Value = '  -11.59%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '8.40'
$ws.Range("E33").Value = '  -11.51%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.87'
$ws.Range("E34").Value = '  -11.68%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.28'
$ws.Range("E35").Value = '  -20.64%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '6.64'
$ws.Range("E36").Value = '  -9.82%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.134'
$ws.Range("E37").Value = '  -9.01%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '22.20'
$ws.Range("E38").Value = '  -5.04%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '21.83'
$ws.Range("E39").Value = '  -0.01%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.998'
$ws.Range("E40").Value = '  -0.14%  '
$ws.Range("E41").Value = '  -2.78%  '
$ws.Range("E42").Value = '  -10.04%  '
$ws.Range("E43").Value = '  -0.02%  '
$ws.Range("E44").Value = '  -12.76%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '147.79'
$ws.Range("E45").Value = '  -7.55%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '174.66'
$ws.Range("E46").Value = '  -9.55%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '43.54'
$ws.Range("E47").Value = '  -6.12%  '
$ws.Range("E48").Value = '  +12.45%  '
$ws.Range("E49").Value = '  -12.19%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '4.03'
$ws.Range("E50").Value = '  -10.20%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.603'
$ws.Range("E51").Value = '  -9.53%  '
